$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers need to be forced
# to remain text (matching the original inlineStr "Price" column formatting)
# by temporarily applying a text NumberFormat before assignment, then restoring
# the default "Normal" style so no stray formatting is left behind.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}


# Row 2
$ws.Range("D2").Value = "43.454.53"
$ws.Range("E2").Value = "  +1.44%  "

# Row 3
$ws.Range("D3").Value = "2.239.26"
$ws.Range("E3").Value = "  +0.55%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
Set-TextValue "D5" "317.59"
$ws.Range("E5").Value = "  +1.47%  "

# Row 6
Set-TextValue "D6" "99.82"
$ws.Range("E6").Value = "  -1.18%  "

# Row 7
$ws.Range("E7").Value = "  +2.43%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
Set-TextValue "D9" "0.564"
$ws.Range("E9").Value = "  +0.06%  "

# Row 10
Set-TextValue "D10" "37.23"
$ws.Range("E10").Value = "  -1.01%  "

# Row 11
Set-TextValue "D11" "0.0831"
$ws.Range("E11").Value = "  -1.73%  "

# Row 12
Set-TextValue "D12" "7.74"
$ws.Range("E12").Value = "  +0.82%  "

# Row 13
$ws.Range("E13").Value = "  +1.68%  "

# Row 14
Set-TextValue "D14" "0.870"
$ws.Range("E14").Value = "  -2.03%  "

# Row 15
$ws.Range("D15").Value = "2.579.69"
$ws.Range("E15").Value = "  +0.82%  "

# Row 16
Set-TextValue "D16" "14.42"
$ws.Range("E16").Value = "  +3.80%  "

# Row 17
$ws.Range("D17").Value = "2.218.28"
$ws.Range("E17").Value = "  +1.06%  "

# Row 18
$ws.Range("D18").Value = "43.397.00"
$ws.Range("E18").Value = "  +1.98%  "

# Row 19
Set-TextValue "D19" "14.18"
$ws.Range("E19").Value = "  -2.22%  "

# Row 20
Set-TextValue "D20" "6.64"
$ws.Range("E20").Value = "  -0.62%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0972"
$ws.Range("E21").Value = "  +1.44%  "

# Row 22
$ws.Range("E22").Value = "  -1.53%  "

# Row 23
Set-TextValue "D23" "65.40"
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
Set-TextValue "D24" "236.87"
$ws.Range("E24").Value = "  +0.11%  "

# Row 25
$ws.Range("E25").Value = "  +1.60%  "

# Row 26
$ws.Range("E26").Value = "  -0.26%  "

# Row 27
$ws.Range("E27").Value = "  +3.19%  "

# Row 28
Set-TextValue "D28" "10.13"
$ws.Range("E28").Value = "  -2.12%  "

# Row 29
$ws.Range("E29").Value = "  +2.88%  "

# Row 30
$ws.Range("E30").Value = "  -3.63%  "

# Row 31
Set-TextValue "D31" "37.02"
$ws.Range("E31").Value = "  +6.72%  "

# Row 32
Set-TextValue "D32" "0.0882"
$ws.Range("E32").Value = "  -2.46%  "

# Row 33
Set-TextValue "D33" "20.35"
$ws.Range("E33").Value = "  -2.04%  "

# Row 34
Set-TextValue "D34" "157.31"
$ws.Range("E34").Value = "  -1.66%  "

# Row 35
Set-TextValue "D35" "2.73"
$ws.Range("E35").Value = "  -0.64%  "

# Row 36
Set-TextValue "D36" "3.25"
$ws.Range("E36").Value = "  +2.39%  "

# Row 37
$ws.Range("E37").Value = "  -0.62%  "

# Row 38
Set-TextValue "D38" "1.89"
$ws.Range("E38").Value = "  -0.82%  "

# Row 39
Set-TextValue "D39" "4.45"
$ws.Range("E39").Value = "  +1.75%  "

# Row 40
$ws.Range("E40").Value = "  +1.44%  "

# Row 41
Set-TextValue "D41" "3.75"
$ws.Range("E41").Value = "  +2.83%  "

# Row 42
Set-TextValue "D42" "0.0323"
$ws.Range("E42").Value = "  +0.12%  "

# Row 43
Set-TextValue "D43" "14.46"
$ws.Range("E43").Value = "  +17.23%  "

# Row 44
$ws.Range("E44").Value = "  +0.16%  "

# Row 45
$ws.Range("D45").Value = "1.797.99"
$ws.Range("E45").Value = "  -0.79%  "

# Row 46
$ws.Range("E46").Value = "  -2.57%  "

# Row 47
Set-TextValue "D47" "85.02"
$ws.Range("E47").Value = "  -7.20%  "

# Row 48
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D48" "5.32"
$ws.Range("E48").Value = "  -1.04%  "

# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D49" "8.89"
$ws.Range("E49").Value = "  +3.89%  "

# Row 50
Set-TextValue "D50" "74.67"
$ws.Range("E50").Value = "  -4.21%  "

# Row 51
Set-TextValue "D51" "58.99"
$ws.Range("E51").Value = "  -4.87%  "
